$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 17 cells: TileMem label + its shared-RAM tile size
$ws.Range("I17").Value = "TileMem"
$ws.Range("J17").Value = 10240

# New formulas in column J (rows 12-14) subtracting the tile reservation
$ws.Range("J12").Formula = '=I12-$J$17'
$ws.Range("J13:J14").Formula = '=I13-$J$17'

# Move the active selection to the newly populated cell
$ws.Range("I17").Select()
